$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly CDC data rows appended below the existing data (rows 80-81)
# Date values are stored as Excel serial date numbers (44443 = 2021-09-04, 44450 = 2021-09-11)
$newRows = @(
    @{ Row = 80; DateSerial = 44443; Values = @(121.1, 228.1, 264.5, 274.2, 218.6, 223.1, 192, 142.9, 109.3, 99.1) },
    @{ Row = 81; DateSerial = 44450; Values = @(89.6, 165.7, 192.3, 198.8, 162.1, 166.1, 145.80000000000001, 109.1, 82.8, 75.5) }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $r.DateSerial
    # Reuse the existing date-formatted style from the row above instead of
    # creating a new number format entry
    $ws.Cells.Item($row - 1, 1).Copy()
    $dateCell.PasteSpecial(-4122)

    $col = 2
    foreach ($val in $r.Values) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
}

$excel.CutCopyMode = $false

# Reflect the updated selection state after the edit
$ws.Range("E75").Select()
